# Update PLC data 2025-10-13 14:16:45
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 180060
$ws.Range("C4").Value = 170011
$ws.Range("C5").Value = 10049
$ws.Range("C8").Value = 65.14
